$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# Rows in column G whose DSL step text still calls the old
# "TakeNativeScreenshot(...)" function. The commit renamed that
# function to "TakeScreenshot(...)" everywhere it was used.
$rows = @(2, 3, 6, 7, 11, 12, 13, 15, 19)

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r, 7)   # column G
    $text = $cell.Value2
    $newText = $text -replace "TakeNativeScreenshot\(", "TakeScreenshot("
    $cell.Value2 = $newText
}

# Update the sheet's active selection from C2 to G2.
$ws.Range("G2").Select()
